$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the list date in A1 (45308 -> 45309, i.e. one day later)
$ws.Range("A1").Value = 45309

# Step 2: update the burlete prices in column D (rows 33-41)
$ws.Range("D33").Value = 698
$ws.Range("D34").Value = 899
$ws.Range("D35").Value = 970
$ws.Range("D36").Value = 1248
$ws.Range("D37").Value = 1578
$ws.Range("D38").Value = 1238
$ws.Range("D39").Value = 1565
$ws.Range("D40").Value = 1961
$ws.Range("D41").Value = 2277
